# Update odds values for the week's matches (Flashscore 2025-06-02 data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Valur - Fram) updates
$ws.Range("H11").Value = 4.85
$ws.Range("I11").Value = 6
$ws.Range("N11").Value = 1.35
$ws.Range("Q11").Value = 3.95
$ws.Range("R11").Value = 1.52
$ws.Range("S11").Value = 2.37
$ws.Range("T11").Value = 12
$ws.Range("U11").Value = 9.75
$ws.Range("W11").Value = 11.25
$ws.Range("AA11").Value = 10.75
$ws.Range("AD11").Value = 28
$ws.Range("AE11").Value = 50
$ws.Range("AF11").Value = 19.5

# Row 12 (Ternana - Pescara) updates
$ws.Range("N12").Value = 2.12
$ws.Range("O12").Value = 1.65
$ws.Range("T12").Value = 6.6
$ws.Range("AB12").Value = 14
